$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 34, shifting existing rows 34:54 down to 35:55.
$ws.Rows("34:34").Insert()

# Populate the newly inserted row 34 with the new record's data.
$ws.Cells.Item(34, 1).Value = 5
$ws.Cells.Item(34, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(34, 3).Value = "Maule"
$ws.Cells.Item(34, 4).Value = 44582
$ws.Cells.Item(34, 5).Value = 7
$ws.Cells.Item(34, 6).Value = "Fruta"
$ws.Cells.Item(34, 7).Value = 100101
$ws.Cells.Item(34, 8).Value = "Berries"
$ws.Cells.Item(34, 9).Value = 100101001
$ws.Cells.Item(34, 10).Value = "Arándano (blue)"
$ws.Cells.Item(34, 11).Value = "Sin especificar"
$ws.Cells.Item(34, 12).Value = "Segunda"
$ws.Cells.Item(34, 13).Value = 160
$ws.Cells.Item(34, 14).Value = 3200
$ws.Cells.Item(34, 15).Value = 3200
$ws.Cells.Item(34, 16).Value = 3200
$ws.Cells.Item(34, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(34, 18).Value = "Provincia de Linares"
$ws.Cells.Item(34, 19).Value = 1600
$ws.Cells.Item(34, 20).Value = 2
